# ランサーズ (sheet1): insert a brand-new scraped row at row 2, pushing all
# existing data rows down by one. Then repair the Hyperlinks collection,
# which does not automatically shift its cell references on row insert.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1) Insert a new row above the current row 2 (the header stays at row 1).
$ws1.Rows.Item(2).Insert()

# 2) Populate the newly inserted row 2 with the new job listing.
$ws1.Range("A2").Value = "2025-08-30 01:40:08"
$ws1.Range("B2").Value = "【急募】kintone業務アプリ開発のプロフェッショナルを探しています!"
$ws1.Range("C2").Value = "システム開発"
$ws1.Range("D2").Value = "300,000 円 ~"
$ws1.Range("E2").Value = "期限情報なし"
$ws1.Range("F2").Value = "https://www.lancers.jp/work/detail/5359873"
$ws1.Range("G2").Value = 93
$ws1.Range("H2").Value = "◆開発 ◇アプリ"

# 3) The engine's row Insert() does not shift the worksheet's Hyperlinks
#    collection along with the data, so rebuild it from scratch: wipe every
#    hyperlink in the sheet (Range.Hyperlinks.Delete() clears the whole
#    sheet) and re-add one per row in the correct, now-shifted order.
$ws1.Range("F2").Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("F2"), "https://www.lancers.jp/work/detail/5359873")
$ws1.Hyperlinks.Add($ws1.Range("F3"), "https://www.lancers.jp/work/detail/5382563")
$ws1.Hyperlinks.Add($ws1.Range("F4"), "https://www.lancers.jp/work/detail/5382589")
$ws1.Hyperlinks.Add($ws1.Range("F5"), "https://www.lancers.jp/work/detail/5382721")
$ws1.Hyperlinks.Add($ws1.Range("F6"), "https://www.lancers.jp/work/detail/5382728")
$ws1.Hyperlinks.Add($ws1.Range("F7"), "https://www.lancers.jp/work/detail/5379313")
$ws1.Hyperlinks.Add($ws1.Range("F8"), "https://www.lancers.jp/work/detail/5382676")
$ws1.Hyperlinks.Add($ws1.Range("F9"), "https://www.lancers.jp/work/detail/5379730")
$ws1.Hyperlinks.Add($ws1.Range("F10"), "https://www.lancers.jp/work/detail/5308620")
$ws1.Hyperlinks.Add($ws1.Range("F11"), "https://www.lancers.jp/work/detail/5309123")
$ws1.Hyperlinks.Add($ws1.Range("F12"), "https://www.lancers.jp/work/detail/5309099")
$ws1.Hyperlinks.Add($ws1.Range("F13"), "https://www.lancers.jp/work/detail/5309115")
$ws1.Hyperlinks.Add($ws1.Range("F14"), "https://www.lancers.jp/work/detail/5309162")
$ws1.Hyperlinks.Add($ws1.Range("F15"), "https://www.lancers.jp/work/detail/5308794")
$ws1.Hyperlinks.Add($ws1.Range("F16"), "https://www.lancers.jp/work/detail/5371807")
$ws1.Hyperlinks.Add($ws1.Range("F17"), "https://www.lancers.jp/work/detail/5309166")
$ws1.Hyperlinks.Add($ws1.Range("F18"), "https://www.lancers.jp/work/detail/5309542")
$ws1.Hyperlinks.Add($ws1.Range("F19"), "https://www.lancers.jp/work/detail/5309147")
$ws1.Hyperlinks.Add($ws1.Range("F20"), "https://www.lancers.jp/work/detail/5308955")
$ws1.Hyperlinks.Add($ws1.Range("F21"), "https://www.lancers.jp/work/detail/5375452")
$ws1.Hyperlinks.Add($ws1.Range("F22"), "https://www.lancers.jp/work/detail/5309229")
$ws1.Hyperlinks.Add($ws1.Range("F23"), "https://www.lancers.jp/work/detail/5308967")
$ws1.Hyperlinks.Add($ws1.Range("F24"), "https://www.lancers.jp/work/detail/5309519")
$ws1.Hyperlinks.Add($ws1.Range("F25"), "https://www.lancers.jp/work/detail/5309439")
$ws1.Hyperlinks.Add($ws1.Range("F26"), "https://www.lancers.jp/work/detail/5309131")
$ws1.Hyperlinks.Add($ws1.Range("F27"), "https://www.lancers.jp/work/detail/5382523")
$ws1.Hyperlinks.Add($ws1.Range("F28"), "https://www.lancers.jp/work/detail/5382629")
$ws1.Hyperlinks.Add($ws1.Range("F29"), "https://www.lancers.jp/work/detail/5382550")
$ws1.Hyperlinks.Add($ws1.Range("F30"), "https://www.lancers.jp/work/detail/5382463")
$ws1.Hyperlinks.Add($ws1.Range("F31"), "https://www.lancers.jp/work/detail/5382289")
$ws1.Hyperlinks.Add($ws1.Range("F32"), "https://www.lancers.jp/work/detail/5371075")
$ws1.Hyperlinks.Add($ws1.Range("F33"), "https://www.lancers.jp/work/detail/5382344")
$ws1.Hyperlinks.Add($ws1.Range("F34"), "https://www.lancers.jp/work/detail/5381977")
$ws1.Hyperlinks.Add($ws1.Range("F35"), "https://www.lancers.jp/work/detail/5382236")
$ws1.Hyperlinks.Add($ws1.Range("F36"), "https://www.lancers.jp/work/detail/5382213")
$ws1.Hyperlinks.Add($ws1.Range("F37"), "https://www.lancers.jp/work/detail/5370186")
$ws1.Hyperlinks.Add($ws1.Range("F38"), "https://www.lancers.jp/work/detail/5381748")
$ws1.Hyperlinks.Add($ws1.Range("F39"), "https://www.lancers.jp/work/detail/5381595")
$ws1.Hyperlinks.Add($ws1.Range("F40"), "https://www.lancers.jp/work/detail/5381608")
$ws1.Hyperlinks.Add($ws1.Range("F41"), "https://www.lancers.jp/work/detail/5379679")
$ws1.Hyperlinks.Add($ws1.Range("F42"), "https://www.lancers.jp/work/detail/5381625")
$ws1.Hyperlinks.Add($ws1.Range("F43"), "https://www.lancers.jp/work/detail/5381634")
$ws1.Hyperlinks.Add($ws1.Range("F44"), "https://www.lancers.jp/work/detail/5381290")
$ws1.Hyperlinks.Add($ws1.Range("F45"), "https://www.lancers.jp/work/detail/5371027")
$ws1.Hyperlinks.Add($ws1.Range("F46"), "https://www.lancers.jp/work/detail/5381284")
$ws1.Hyperlinks.Add($ws1.Range("F47"), "https://www.lancers.jp/work/detail/5381264")
$ws1.Hyperlinks.Add($ws1.Range("F48"), "https://www.lancers.jp/work/detail/5381250")
$ws1.Hyperlinks.Add($ws1.Range("F49"), "https://www.lancers.jp/work/detail/5381245")
$ws1.Hyperlinks.Add($ws1.Range("F50"), "https://www.lancers.jp/work/detail/5381204")
$ws1.Hyperlinks.Add($ws1.Range("F51"), "https://www.lancers.jp/work/detail/5381118")
$ws1.Hyperlinks.Add($ws1.Range("F52"), "https://www.lancers.jp/work/detail/5380896")
$ws1.Hyperlinks.Add($ws1.Range("F53"), "https://www.lancers.jp/work/detail/5380830")
$ws1.Hyperlinks.Add($ws1.Range("F54"), "https://www.lancers.jp/work/detail/5016989")
$ws1.Hyperlinks.Add($ws1.Range("F55"), "https://www.lancers.jp/work/detail/5273634")
$ws1.Hyperlinks.Add($ws1.Range("F56"), "https://www.lancers.jp/work/detail/5217096")
$ws1.Hyperlinks.Add($ws1.Range("F57"), "https://www.lancers.jp/work/detail/5380343")
$ws1.Hyperlinks.Add($ws1.Range("F58"), "https://www.lancers.jp/work/detail/5380337")
$ws1.Hyperlinks.Add($ws1.Range("F59"), "https://www.lancers.jp/work/detail/5380683")
$ws1.Hyperlinks.Add($ws1.Range("F60"), "https://www.lancers.jp/work/detail/5380747")
$ws1.Hyperlinks.Add($ws1.Range("F61"), "https://www.lancers.jp/work/detail/5341051")
$ws1.Hyperlinks.Add($ws1.Range("F62"), "https://www.lancers.jp/work/detail/5380357")
$ws1.Hyperlinks.Add($ws1.Range("F63"), "https://www.lancers.jp/work/detail/5380420")

# 統計 (sheet2): append one new summary row at the bottom (row 24).
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A24").Value = "2025-08-30T01:40:08.714767"
$ws2.Range("B24").Value = 24
$ws2.Range("C24").Value = "全案件リスト"
$ws2.Range("D24").Value = 62.5
$ws2.Range("E24").Value = 9
$ws2.Range("F24").Value = 10
$ws2.Range("G24").Value = 24
